$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column E: "Date Created (Year)*" header with two data rows of 2000,
# plus an extra (empty) formatted row below, matching the black explicit font color style.
$ws.Range("E1:E4").Font.Color = 0
$ws.Range("E1").Value = "Date Created (Year)*"
$ws.Range("E2").Value = 2000
$ws.Range("E3").Value = 2000

# Update the selection/view: select D1 (this also clears the old C1 scroll position).
$ws.Range("D1").Select()
